# The document carries a set of SharePoint "documentManagement" custom XML
# parts (content-type schema, form templates, and their associated
# datastore/property parts). These were only ever injected by the
# SharePoint document library the .docx was checked out from - they are not
# referenced anywhere in the visible document content - so remove them from
# the package entirely, walking backwards through the collection (the
# standard COM pattern) so deleting an item doesn't disturb the indices of
# the ones still to be visited.
$d = $word.ActiveDocument

$customXmlParts = $d.CustomXMLParts
for ($i = $customXmlParts.Count; $i -ge 1; $i--) {
    $customXmlParts.Item($i).Delete()
}

Write-Output ("CustomXMLParts remaining: " + $d.CustomXMLParts.Count)
